$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GASS")

# Insert a new column before column D, shifting existing D:K data to E:L.
$ws.Range("D1").EntireColumn.Insert(-4161, 0)

# The newly inserted column D is blank with default formatting; copy the
# number formats/styles from column E (the original column D, now shifted)
# so the new quarter column matches the look of its neighbors.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# The format-only paste above can stamp a phantom blank cell into column D
# on rows that have no data at all (section header / spacer rows). Remove
# those so the row layout matches the source rows exactly.
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate the new (most recent) quarter column D with its reported values.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 42700
$ws.Range("D9").Value = 21400
$ws.Range("D10").Value = 21300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 600
$ws.Range("D15").Value = 10100
$ws.Range("D17").Value = 37600
$ws.Range("D18").Value = 5100
$ws.Range("D20").Value = 200
$ws.Range("D21").Value = 15500
$ws.Range("D22").Value = 6100
$ws.Range("D23").Value = -800
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -800
$ws.Range("D27").Value = -800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -200
$ws.Range("D33").Value = -800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -800
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 64800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 4500
$ws.Range("D44").Value = 4200
$ws.Range("D45").Value = 17100
$ws.Range("D46").Value = 90600
$ws.Range("D47").Value = 100
$ws.Range("D48").Value = 953300
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 14800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1058800
$ws.Range("D57").Value = 22600
$ws.Range("D58").Value = 48600
$ws.Range("D59").Value = 13400
$ws.Range("D60").Value = 84600
$ws.Range("D61").Value = 406100
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 490800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 86200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 568000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -800
$ws.Range("D83").Value = 10100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 9500
$ws.Range("D91").Value = -200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 17700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -18100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 9100
